$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (it spans from the
#    "Accomplishments at 1st code sprint" heading down to the end of the
#    "Discussed order of calls..." paragraph).  We delete it now - before
#    inserting the new paragraph below - so that the name lookup
#    unambiguously finds this (the only) "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# ---------------------------------------------------------------------------
# 2) Insert a new bulleted paragraph right after the "Documentation" bullet
#    (under the "Misc:" heading), containing the new note, and re-create the
#    "_GoBack" bookmark (now empty / collapsed) at the end of its text.
# ---------------------------------------------------------------------------
$docParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Documentation`r") {
        [void]($docParaIndex = $i)
        break
    }
}
Write-Output "Documentation paragraph index: $docParaIndex"

$docPara = $d.Paragraphs.Item($docParaIndex)
$insertRange = $docPara.Range
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)
$newPara = $d.Paragraphs.Item($docParaIndex + 1)

$newParaXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Compare results from the original version of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sagehen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (i.e., GSFLOW-only, no lakes) using the published GSFLOW code and the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Fortran</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> version of MMF GSFLOW code.  This is to make sure that </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$newPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 3) Move the "lastRenderedPageBreak" marker: drop it from the run that
#    starts "Got started with repos on Monday" and add it to the run that
#    starts the "Accomplishments at 2nd code sprint..." heading.
# ---------------------------------------------------------------------------
$accompRange = $d.Content
$accompRange.Find.Execute("Accomplishments at 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$accompRange.Delete()
$accompXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00181968"><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Accomplishments at 2</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$accompInsertPoint = $d.Range($accompRange.Start, $accompRange.Start)
$accompInsertPoint.InsertXML($accompXml)

$mondayRange = $d.Content
$mondayRange.Find.Execute("Got started with repos on Monday", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mondayRange.Delete()
$mondayXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Got started with repos on Monday</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$mondayInsertPoint = $d.Range($mondayRange.Start, $mondayRange.Start)
$mondayInsertPoint.InsertXML($mondayXml)

Write-Output "Edit complete"
